# Applies "first correction SFN snippet" changes to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variable Snippets")

# G3: SFN_stamp.xml -> SFNSyslogServer_stamp.xml
$ws.Range("G3").Value = "SFNSyslogServer_stamp.xml"

# Column E repeating block (rows 6,11,16,...,81): fix placeholder spacing
$startRows = 6,11,16,21,26,31,36,41,46,51,56,61,66,71,76,81

foreach ($r in $startRows) {
    $ws.Cells.Item($r, 5).Value = "{{ ObjName }}"
    $ws.Cells.Item($r + 1, 5).Value = "{{ ObjIPNetmask }}"
    $ws.Cells.Item($r + 2, 5).Value = "{{ ObjDescription }}"
    $ws.Cells.Item($r + 3, 5).Value = "{{ TagName }}"
}

# Update the pane's active cell selection (cosmetic) to G8
$ws.Range("G8").Select()

# Update workbook window position (cosmetic)
$excel.Left = 135
$excel.Top = 450
